# Apply the recorded edits to testdata.xlsx:
#  - RUNMANAGER sheet: extend the selection from A2 to A2:E5 (active cell stays A2)
#  - DATA sheet: move the selection/active cell from C7 to B3
#  - DATA sheet: change cells B4 and B6 from "yes" to "no"

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# RUNMANAGER: select A2:E5 (this naturally makes A2 the active cell)
$wsRunManager.Activate()
$wsRunManager.Range("A2:E5").Select()

# DATA: update values and selection; leave DATA as the active/selected sheet
$wsData.Activate()
$wsData.Range("B4").Value = "no"
$wsData.Range("B6").Value = "no"
$wsData.Range("B3").Select()
